# Adding Multi Browser execution feature to Framework
# Insert a new "browser" column into the DATA sheet (between the
# "execute" and "username" columns), populate it with the browsers used
# by each test row, and flip the second test row back to "yes" so it
# also executes (it previously opted out via "No").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")
$ws.Activate()

# Insert a new column at C, shifting username/password right.
$ws.Columns("C").Insert()

# Header for the new column.
$ws.Range("C1").Value = "browser"

# Per-row browser values.
$ws.Range("C2").Value = "chrome"
$ws.Range("C3").Value = "edge"

# Row 3 now also executes (was "No").
$ws.Range("B3").Value = "yes"

# Match the selection left behind in the authored workbook.
$ws.Range("A4:F18").Select()
